$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be parsed as a number: force text via NumberFormat, then restore default style
$forceCells = @("D5","D7","D9","D10","D11","D12","D14","D15","D20","D21","D22","D24","D27","D29","D32","D33","D34","D35","D36","D37","D40","D41","D42","D43","D44","D45","D47","D50","D51")
foreach ($addr in $forceCells) { $ws.Range($addr).NumberFormat = "@" }

$ws.Range('D5').Value2 = '234.04'
$ws.Range('D7').Value2 = '71.14'
$ws.Range('D9').Value2 = '0.491'
$ws.Range('D10').Value2 = '0.0973'
$ws.Range('D11').Value2 = '27.20'
$ws.Range('D12').Value2 = '0.106'
$ws.Range('D14').Value2 = '16.07'
$ws.Range('D15').Value2 = '6.27'
$ws.Range('D20').Value2 = '6.32'
$ws.Range('D21').Value2 = '74.32'
$ws.Range('D22').Value2 = '249.53'
$ws.Range('D24').Value2 = '3.77'
$ws.Range('D27').Value2 = '10.08'
$ws.Range('D29').Value2 = '172.37'
$ws.Range('D32').Value2 = '0.126'
$ws.Range('D33').Value2 = '4.97'
$ws.Range('D34').Value2 = '0.0688'
$ws.Range('D35').Value2 = '5.05'
$ws.Range('D36').Value2 = '3.69'
$ws.Range('D37').Value2 = '6.53'
$ws.Range('D40').Value2 = '1.00'
$ws.Range('D41').Value2 = '18.69'
$ws.Range('D42').Value2 = '8.88'
$ws.Range('D43').Value2 = '4.48'
$ws.Range('D44').Value2 = '1.16'
$ws.Range('D45').Value2 = '98.87'
$ws.Range('D47').Value2 = '0.0959'
$ws.Range('D50').Value2 = '2.76'
$ws.Range('D51').Value2 = '0.000203'

foreach ($addr in $forceCells) { $ws.Range($addr).Style = "Normal" }

# Remaining cells are already non-numeric text (URLs, names, padded percentages) -- plain assignment keeps them as text
$ws.Range('D2').Value2 = '43.219.68'
$ws.Range('E2').Value2 = '  -0.92%  '
$ws.Range('D3').Value2 = '2.355.62'
$ws.Range('E3').Value2 = '  +4.99%  '
$ws.Range('E4').Value2 = '  -0.03%  '
$ws.Range('E5').Value2 = '  +1.44%  '
$ws.Range('E6').Value2 = '  +0.20%  '
$ws.Range('E7').Value2 = '  +12.57%  '
$ws.Range('E8').Value2 = '  +0.04%  '
$ws.Range('E9').Value2 = '  +11.34%  '
$ws.Range('E10').Value2 = '  +1.02%  '
$ws.Range('E11').Value2 = '  +2.34%  '
$ws.Range('B12').Value2 = 'TRON'
$ws.Range('C12').Value2 = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('E12').Value2 = '  +1.52%  '
$ws.Range('B13').Value2 = 'WrappedliquidstakedEther2.0'
$ws.Range('C13').Value2 = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D13').Value2 = '2.707.09'
$ws.Range('E13').Value2 = '  +4.86%  '
$ws.Range('E14').Value2 = '  +3.71%  '
$ws.Range('E15').Value2 = '  +3.77%  '
$ws.Range('E16').Value2 = '  +3.84%  '
$ws.Range('D17').Value2 = '2.349.57'
$ws.Range('E17').Value2 = '  +4.41%  '
$ws.Range('D18').Value2 = '43.251.50'
$ws.Range('E18').Value2 = '  -0.65%  '
$ws.Range('E19').Value2 = '  +3.87%  '
$ws.Range('E20').Value2 = '  +3.70%  '
$ws.Range('E21').Value2 = '  +1.51%  '
$ws.Range('E22').Value2 = '  +1.42%  '
$ws.Range('E23').Value2 = '  +0.02%  '
$ws.Range('E24').Value2 = '  +1.36%  '
$ws.Range('E25').Value2 = '  +1.45%  '
$ws.Range('E26').Value2 = '  -1.51%  '
$ws.Range('E27').Value2 = '  +2.90%  '
$ws.Range('E28').Value2 = '  +3.64%  '
$ws.Range('E29').Value2 = '  -0.21%  '
$ws.Range('E30').Value2 = '  +7.60%  '
$ws.Range('E31').Value2 = '  -2.84%  '
$ws.Range('E32').Value2 = '  +0.19%  '
$ws.Range('E33').Value2 = '  +2.22%  '
$ws.Range('E34').Value2 = '  +2.08%  '
$ws.Range('E35').Value2 = '  +3.50%  '
$ws.Range('E36').Value2 = '  +1.66%  '
$ws.Range('E37').Value2 = '  +3.20%  '
$ws.Range('E38').Value2 = '  +6.60%  '
$ws.Range('E39').Value2 = '  +0.77%  '
$ws.Range('B40').Value2 = 'BinanceUSD'
$ws.Range('C40').Value2 = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('E40').Value2 = '  -0.13%  '
$ws.Range('B41').Value2 = 'InjectiveProtocol'
$ws.Range('C41').Value2 = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('E41').Value2 = '  +9.50%  '
$ws.Range('E42').Value2 = '  +3.76%  '
$ws.Range('B43').Value2 = 'FTXToken'
$ws.Range('C43').Value2 = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('E43').Value2 = '  -1.16%  '
$ws.Range('B44').Value2 = 'ARBITRUM'
$ws.Range('C44').Value2 = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('E44').Value2 = '  +8.58%  '
$ws.Range('E45').Value2 = '  +1.29%  '
$ws.Range('E46').Value2 = '  +2.19%  '
$ws.Range('E47').Value2 = '  +2.71%  '
$ws.Range('D48').Value2 = '1.438.59'
$ws.Range('E48').Value2 = '  -0.69%  '
$ws.Range('D49').Value2 = '2.579.37'
$ws.Range('E49').Value2 = '  +5.11%  '
$ws.Range('E50').Value2 = '  +0.24%  '
$ws.Range('E51').Value2 = '  -2.48%  '
